# =========================================================================
# LOT2056.xlsx restructuring
#
# The "Docentes responsaveis" / "Programa" / "Avaliacao" block got reshuffled:
#   - the teacher name that used to sit in its own unlabeled row now fills in
#     the "Objetivos:" row, and resurfaces again under "Programa:"/"Programa
#     resumido:" (reusing the same text the source workbook already had).
#   - "Avaliacao:" / "Requisitos:" become single-column label-only rows.
#   - the two trailing "Requisitos" detail rows (24 and 25) are removed, so the
#     sheet dimension shrinks from A1:C25 to A1:C23.
# =========================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Cells that keep their (row, col) position but change content -----
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range("B2").Value = 'LOT2056'
$ws.Range("C2").Value = 'LOT2056'
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Trabalho de Conclusão de Curso I'
$ws.Range("C3").Value = ' Trabalho de Conclusão de Curso I'
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Course Completion Work I'
$ws.Range("C4").Value = 'Course Completion Work I'
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '45 h'
$ws.Range("C7").Value = '45 h'
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EB-9'
$ws.Range("C9").Value = 'EB-9'
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C10").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'Lead students to develop a course completion project on a specific topic related to Biochemical engineering.'
$ws.Range("C11").Value = 'Lead students to develop a course completion project on a specific topic related to Biochemical engineering.'
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B14").Value = '1) Research Methodology in Engineering. 2) Elements of a Research Project. 3) Research Methods. 4) Guidelines of a research project.5) Steps in writing scientific texts'
$ws.Range("C14").Value = '1) Research Methodology in Engineering. 2) Elements of a Research Project. 3) Research Methods. 4) Guidelines of a research project.5) Steps in writing scientific texts'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C15").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. Research Methodology in Engineering: principles and definition.2. Scientific Project: what is a research project. The elements that make up a research project.3. Research Methods used in Biochemical Engineering.4. Rules for preparation of text and Bibliographical References.5. Steps in writing scientific texts'
$ws.Range("C16").Value = '1. Research Methodology in Engineering: principles and definition.2. Scientific Project: what is a research project. The elements that make up a research project.3. Research Methods used in Biochemical Engineering.4. Rules for preparation of text and Bibliographical References.5. Steps in writing scientific texts'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '8853480 - Tatiane da Franca Silva'
$ws.Range("C18").Value = '8853480 - Tatiane da Franca Silva'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'M=≥ 5,0 para ser aprovado'
$ws.Range("C20").Value = 'M=≥ 5,0 para ser aprovado'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = '(NF+RP)/2 ≥ 5,0 para ser aprovado, onde RP é a nota do projeto modificado apresentado.'
$ws.Range("C21").Value = '(NF+RP)/2 ≥ 5,0 para ser aprovado, onde RP é a nota do projeto modificado apresentado.'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B23").Value = 'LOT2013 -  Engenharia Bioquímica I  (Requisito fraco)' + [char]10 + ''
$ws.Range("C23").Value = 'LOT2013 -  Engenharia Bioquímica I  (Requisito fraco)' + [char]10 + ''

# Cells whose new text would otherwise be auto-detected as a number/date by
# Excel (e.g. "1" or "01/01/2020"); force them to remain plain text by (a)
# pre-formatting as Text, (b) assigning the value, then (c) pasting the
# number-format/style back from an untouched neighbor in the same column so
# the final style index matches the rest of the column.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '1'
$ws.Range("B7").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '1'
$ws.Range("C7").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '1'
$ws.Range("B7").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '1'
$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '01/01/2020'
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '01/01/2020'
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '01/01/2020'
$ws.Range("B7").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '01/01/2020'
$ws.Range("C7").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- 2) Cells that must disappear entirely (row shrinks to fewer columns) -
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# ---- 3) Brand-new cells that did not exist at this grid position before ---
# Column A picks up the correct style (s=1) automatically from the column
# default. The new B19/C19 cells need their format copied from the already-
# correctly-styled B18/C18 cells so they end up with s=2/s=3 like the rest
# of their columns.
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B19").Value = 'Apresentação de um pré-projeto e um projeto. O projeto será avaliado por dois examinadores. A média (M) será calculada levando-se a nota do pré-projeto(NPP) e a média da nova dos dois examinadores (NP) M = (0,3PP + 0,7NP), conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica.'
$ws.Range("C19").Value = 'Apresentação de um pré-projeto e um projeto. O projeto será avaliado por dois examinadores. A média (M) será calculada levando-se a nota do pré-projeto(NPP) e a média da nova dos dois examinadores (NP) M = (0,3PP + 0,7NP), conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica.'
$ws.Range("B18:C18").Copy()
$ws.Range("B19:C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- 4) Row heights for the new content layout -----------------------------
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30

# Rows 17 and 22 lose their custom height now that they only hold a single
# label cell in column A; AutoFit drops back to the sheet default (15pt) and
# clears the customHeight flag entirely.
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()

# ---- 5) Drop the two trailing rows so dimension becomes A1:C23 ------------
$ws.Range("A24:C25").EntireRow.Delete()

$ws.Range("A1").Select()
